$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q1" sheet (fund-holdings detail) so the
#    old data is preserved under its original name, then repurpose the
#    original sheet (keeps the same sheetId/rId) as the new "2022-Q4"
#    sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Copy($null, $q1)

$q1Copy = $wb.Worksheets.Item("2022-Q1 (2)")
$q4 = $wb.Worksheets.Item("2022-Q1")

$q4.Name = "2022-Q4"
$q1Copy.Name = "2022-Q1"

$total = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 2) Trim the "2022-Q4" sheet down to just a header row + 2 data rows
#    (it still has the 5-row Q1 detail at this point), then overwrite
#    with the Q4 fund-holdings data.
# ---------------------------------------------------------------------
$q4.Rows.Item(5).Delete()
$q4.Rows.Item(4).Delete()
$q4.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# 3) Re-style the "2022-Q4" sheet's header row / index column to match
#    the bold+border style used elsewhere for section headers (same
#    style already used on the "总计" sheet's header row).
# ---------------------------------------------------------------------
$total.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A3").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "015553"
$q4.Range("B2").ClearFormats()
$q4.Range("C2").Value = "融通价值成长混合A"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "3.54"
$q4.Range("D2").ClearFormats()
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "94.64"
$q4.Range("E2").ClearFormats()
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "2.94"
$q4.Range("F2").ClearFormats()
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.1041"
$q4.Range("G2").ClearFormats()
$q4.Range("H2").Value = 9

$q4.Range("A3").Value = 1
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "015554"
$q4.Range("B3").ClearFormats()
$q4.Range("C3").Value = "融通价值成长混合C"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "1.36"
$q4.Range("D3").ClearFormats()
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "94.64"
$q4.Range("E3").ClearFormats()
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "2.94"
$q4.Range("F3").ClearFormats()
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0400"
$q4.Range("G3").ClearFormats()
$q4.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 4) Update the "总计" (summary) sheet: shift the existing 2022-Q1 row
#    down to row 3, and add the new 2022-Q4 totals in row 2.
# ---------------------------------------------------------------------
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.03

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.14
